# Update excel for forecast
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (row 6 and row 7)
$ws.Range("A6").Value = 45042
$ws.Range("B6").Value = 59151
$ws.Range("C6").Value = "800 urls"

$ws.Range("A7").Value = 45043
$ws.Range("B7").Value = 60269
$ws.Range("C7").Value = "800 urls"

# Match the date cell formatting used by the existing rows above
$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the forecast formula in B25 to include the new rows
$ws.Range("B25").Formula = "=_xlfn.FORECAST.LINEAR(A25,B1:B7,A1:A7)"

# Update the selected cell
$ws.Range("B11").Select()
